$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3578975.5
$ws.Range("J17").Value = 3711459.8
$ws.Range("L17").Value = 11134379.4
$ws.Range("N17").Value = -11134715.4
$ws.Range("H41").Value = 374.8889
$ws.Range("J41").Value = 514
$ws.Range("L41").Value = 514
$ws.Range("N41").Value = -1394
$ws.Range("I98").Value = 623.75
$ws.Range("K98").Value = 623.75
$ws.Range("M98").Value = 874.25
$ws.Range("H112").Value = 4445487.5
$ws.Range("I112").Value = 770
$ws.Range("K112").Value = 2310
$ws.Range("M112").Value = -1202
$ws.Range("I122").Value = 623.75
$ws.Range("K122").Value = 1871.25
$ws.Range("M122").Value = 578.75
$ws.Range("H124").Value = 30924.143
$ws.Range("J124").Value = 31293.334
$ws.Range("L124").Value = 31293.334
$ws.Range("N124").Value = -41113.334
$ws.Range("H125").Value = 1218
$ws.Range("I125").Value = 400
$ws.Range("K125").Value = 3600
$ws.Range("M125").Value = -1140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2895.0833
$ws.Range("I45").Value = 2559.8
$ws.Range("K45").Value = 2559.8
$ws.Range("M45").Value = -2182.8
$ws.Range("H61").Value = 5765
$ws.Range("I61").Value = 6122.222
$ws.Range("K61").Value = 6122.222
$ws.Range("M61").Value = -5910.222
$ws.Range("H122").Value = 2131.4783
$ws.Range("I122").Value = 1264.5264
$ws.Range("J122").Value = 6249.5
$ws.Range("K122").Value = 3793.5792
$ws.Range("L122").Value = 18748.5
$ws.Range("M122").Value = -1343.5792
$ws.Range("N122").Value = -23648.5
$ws.Range("H136").Value = 5765
$ws.Range("I136").Value = 6122.222
$ws.Range("K136").Value = 18366.666
$ws.Range("M136").Value = -15816.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4049.5356
$ws.Range("I134").Value = 4155.074
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 12465.222
$ws.Range("L134").Value = 3600
$ws.Range("M134").Value = -9930.221999999998
$ws.Range("N134").Value = -8670

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17730.064
$ws.Range("I58").Value = 1428.8235
$ws.Range("J58").Value = 37524.43
$ws.Range("K58").Value = 1428.8235
$ws.Range("L58").Value = 37524.43
$ws.Range("M58").Value = -1225.8235
$ws.Range("N58").Value = -37930.43
$ws.Range("H122").Value = 1695.1666
$ws.Range("I122").Value = 1890.3334
$ws.Range("K122").Value = 5671.0002
$ws.Range("M122").Value = -3221.0002
$ws.Range("H134").Value = 1359.8572
$ws.Range("I134").Value = 1265.9231
$ws.Range("J134").Value = 1512.5
$ws.Range("K134").Value = 3797.7693
$ws.Range("L134").Value = 4537.5
$ws.Range("M134").Value = -1262.7693
$ws.Range("N134").Value = -9607.5
$ws.Range("H136").Value = 17730.064
$ws.Range("I136").Value = 1428.8235
$ws.Range("J136").Value = 37524.43
$ws.Range("K136").Value = 4286.470499999999
$ws.Range("L136").Value = 112573.29
$ws.Range("M136").Value = -1736.470499999999
$ws.Range("N136").Value = -117673.29

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 3000
$ws.Range("M14").Value = -2827
$ws.Range("H23").Value = 411.66666
$ws.Range("I23").Value = 40
$ws.Range("J23").Value = 597.5
$ws.Range("K23").Value = 120
$ws.Range("L23").Value = 1792.5
$ws.Range("M23").Value = 115
$ws.Range("N23").Value = -2262.5
$ws.Range("H32").Value = 2000
$ws.Range("J32").Value = 2000
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6566
$ws.Range("H33").Value = 66.833336
$ws.Range("J33").Value = 93.75
$ws.Range("L33").Value = 562.5
$ws.Range("N33").Value = -1128.5
$ws.Range("H120").Value = 17507.5
$ws.Range("I120").Value = 10030
$ws.Range("K120").Value = 30090
$ws.Range("M120").Value = -25252
$ws.Range("H123").Value = 6985
$ws.Range("J123").Value = 6985
$ws.Range("L123").Value = 20955
$ws.Range("N123").Value = -25855
$ws.Range("H131").Value = 710.96
$ws.Range("I131").Value = 546
$ws.Range("J131").Value = 719.6421
$ws.Range("K131").Value = 1638
$ws.Range("L131").Value = 2158.9263
$ws.Range("M131").Value = 3402
$ws.Range("N131").Value = -12238.9263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 12575.167
$ws.Range("J93").Value = 12575.167
$ws.Range("L93").Value = 12575.167
$ws.Range("N93").Value = -16319.167
$ws.Range("H122").Value = 2492.5
$ws.Range("I122").Value = 2390.2
$ws.Range("J122").Value = 3004
$ws.Range("K122").Value = 7170.599999999999
$ws.Range("L122").Value = 9012
$ws.Range("M122").Value = -4720.599999999999
$ws.Range("N122").Value = -13912

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7457
$ws.Range("I7").Value = 4099.8
$ws.Range("K7").Value = 4099.8
$ws.Range("M7").Value = -3987.8
$ws.Range("H36").Value = 30715
$ws.Range("J36").Value = 30715
$ws.Range("L36").Value = 30715
$ws.Range("N36").Value = -31839
$ws.Range("H40").Value = 3300.0344
$ws.Range("I40").Value = 2626.348
$ws.Range("J40").Value = 5882.5
$ws.Range("K40").Value = 2626.348
$ws.Range("L40").Value = 5882.5
$ws.Range("M40").Value = -2490.348
$ws.Range("N40").Value = -6154.5
$ws.Range("H122").Value = 936529.9
$ws.Range("I122").Value = 2803849.2
$ws.Range("J122").Value = 2870.2856
$ws.Range("K122").Value = 8411547.600000001
$ws.Range("L122").Value = 8610.856800000001
$ws.Range("M122").Value = -8409097.600000001
$ws.Range("N122").Value = -13510.8568
$ws.Range("H126").Value = 7457
$ws.Range("I126").Value = 4099.8
$ws.Range("K126").Value = 12299.4
$ws.Range("M126").Value = -9829.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1205.9166
$ws.Range("I126").Value = 1247.238
$ws.Range("K126").Value = 3741.714
$ws.Range("M126").Value = -1271.714
$ws.Range("H132").Value = 1798.9131
$ws.Range("I132").Value = 1054.8125
$ws.Range("J132").Value = 3499.7144
$ws.Range("K132").Value = 3164.4375
$ws.Range("L132").Value = 10499.1432
$ws.Range("M132").Value = -634.4375
$ws.Range("N132").Value = -15559.1432
